$d = $word.ActiveDocument

# Locate the paragraph that marks the end of the content we keep
# ("Restart your browser and you should be able to see the applets that
# use Java 3D.") and the paragraph that marks the end of the block we
# want removed ("Click Continue on the Security Warning dialog").
$keepEndText   = "Restart your browser and you should be able to see the applets that use Java 3D."
$removeEndText = "Click Continue on the Security Warning dialog"

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $keepEndText) {
        $startPara = $i
    }
    if ($t -eq $removeEndText) {
        $endPara = $i
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate the anchor paragraphs for the deletion range."
}

# Range to remove: from the start of the paragraph right after the
# "Restart your browser..." paragraph, through the end (including the
# paragraph mark) of the "Click Continue..." paragraph.
$rangeStart = $d.Paragraphs.Item($startPara + 1).Range.Start
$rangeEnd   = $d.Paragraphs.Item($endPara).Range.End

$r = $d.Range($rangeStart, $rangeEnd)
$r.Delete()
